$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BICT")
$ws.Range("A1").Value = "test"
